$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Totali" summary table headers (row 19) ---
$ws.Range("D19").Value2 = "Totale fallimenti"
$ws.Range("E19").Value2 = "Generazioni mancanti"
$ws.Range("F19").Value2 = "Generazioni non necessarie"
$ws.Range("G19").Value2 = "Generazioni necessarie ma errate"

# --- Fill in new data for the "LLM" row (row 20) ---
$ws.Range("E20").Value2 = 1
$ws.Range("F20").Value2 = 0
$ws.Range("G20").Value2 = 0

# --- Fill in new data for the "Analitica" row (row 21) ---
$ws.Range("E21").Value2 = 1
$ws.Range("F21").Value2 = 0
$ws.Range("G21").Value2 = 0

# --- Widen column G to fit the new, longer header text ---
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(7).ColumnWidth + 6

# --- Move the selection / scrolled view to G21 ---
$ws.Range("G21").Select()
